$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.123382
$ws.Range("N2").Value = 0.246764
$ws.Range("Q2").Value = 0.1496800507533333
$ws.Range("R2").Value = 0.89808030452
